$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 9999
$ws.Range("I32").Value = 9999
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 9999
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -9673
$ws.Range("N32").ClearContents()
$ws.Range("H51").Value = 4365.44
$ws.Range("I51").Value = 2806.8572
$ws.Range("K51").Value = 2806.8572
$ws.Range("M51").Value = -2322.8572
$ws.Range("H69").Value = 5248.75
$ws.Range("H72").Value = 5248.75
$ws.Range("H95").Value = 38127.4
$ws.Range("J95").Value = 38127.4
$ws.Range("L95").Value = 38127.4
$ws.Range("N95").Value = -43619.4
$ws.Range("H138").Value = 5742.691
$ws.Range("J138").Value = 7163.7573
$ws.Range("L138").Value = 21491.2719
$ws.Range("N138").Value = -31771.2719
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 15153240
$ws.Range("I61").Value = 17858714
$ws.Range("J61").Value = 2582.2
$ws.Range("K61").Value = 17858714
$ws.Range("L61").Value = 2582.2
$ws.Range("M61").Value = -17858502
$ws.Range("N61").Value = -3006.2
$ws.Range("H109").Value = 60188.5
$ws.Range("J109").Value = 60188.5
$ws.Range("L109").Value = 60188.5
$ws.Range("N109").Value = -62962.5
$ws.Range("H110").Value = 16047.107
$ws.Range("I110").Value = 16511.385
$ws.Range("J110").Value = 10011.5
$ws.Range("K110").Value = 16511.385
$ws.Range("L110").Value = 10011.5
$ws.Range("M110").Value = -14466.385
$ws.Range("N110").Value = -14101.5
$ws.Range("H132").Value = 31304566
$ws.Range("I132").Value = 13621.333
$ws.Range("K132").Value = 40863.999
$ws.Range("M132").Value = -38333.999
$ws.Range("H136").Value = 15153240
$ws.Range("I136").Value = 17858714
$ws.Range("J136").Value = 2582.2
$ws.Range("K136").Value = 53576142
$ws.Range("L136").Value = 7746.599999999999
$ws.Range("M136").Value = -53573592
$ws.Range("N136").Value = -12846.6
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H38").Value = 20163.5
$ws.Range("I38").Value = 332
$ws.Range("K38").Value = 332
$ws.Range("M38").Value = 84
$ws.Range("H86").Value = 10578.462
$ws.Range("I86").Value = 11865.546
$ws.Range("K86").Value = 11865.546
$ws.Range("M86").Value = -10742.546
$ws.Range("H89").Value = 10578.462
$ws.Range("I89").Value = 11865.546
$ws.Range("K89").Value = 59327.73
$ws.Range("M89").Value = -53711.73
$ws.Range("H108").Value = 60000
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 60000
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 60000
$ws.Range("M108").ClearContents()
$ws.Range("N108").Value = -67680
$ws.Range("H134").Value = 2295.625
$ws.Range("I134").Value = 1978.6842
$ws.Range("J134").Value = 3500
$ws.Range("K134").Value = 5936.0526
$ws.Range("L134").Value = 10500
$ws.Range("M134").Value = -3401.0526
$ws.Range("N134").Value = -15570
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 9693.634
$ws.Range("I99").Value = 6542.6665
$ws.Range("J99").Value = 11794.277
$ws.Range("K99").Value = 6542.6665
$ws.Range("L99").Value = 11794.277
$ws.Range("M99").Value = -5044.6665
$ws.Range("N99").Value = -14790.277
$ws.Range("H103").Value = 15807.667
$ws.Range("I103").Value = 15807.667
$ws.Range("K103").Value = 15807.667
$ws.Range("M103").Value = -14635.667
$ws.Range("H105").Value = 7812.64
$ws.Range("I105").Value = 3696.4092
$ws.Range("K105").Value = 3696.4092
$ws.Range("M105").Value = -1949.4092
$ws.Range("H109").Value = 46249.875
$ws.Range("J109").Value = 43333.332
$ws.Range("L109").Value = 43333.332
$ws.Range("N109").Value = -45413.332
$ws.Range("H122").Value = 5849915
$ws.Range("I122").Value = 2153
$ws.Range("J122").Value = 10528125
$ws.Range("K122").Value = 6459
$ws.Range("L122").Value = 31584375
$ws.Range("M122").Value = -4009
$ws.Range("N122").Value = -31589275
$ws.Range("H126").Value = 9693.634
$ws.Range("I126").Value = 6542.6665
$ws.Range("J126").Value = 11794.277
$ws.Range("K126").Value = 19627.9995
$ws.Range("L126").Value = 35382.831
$ws.Range("M126").Value = -17157.9995
$ws.Range("N126").Value = -40322.831
$ws.Range("H131").Value = 32600
$ws.Range("J131").Value = 32600
$ws.Range("L131").Value = 32600
$ws.Range("N131").Value = -42680
$ws.Range("H134").Value = 3574.4531
$ws.Range("I134").Value = 3649.0847
$ws.Range("J134").Value = 2693.8
$ws.Range("K134").Value = 10947.2541
$ws.Range("L134").Value = 8081.400000000001
$ws.Range("M134").Value = -8412.2541
$ws.Range("N134").Value = -13151.4
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1068.4897
$ws.Range("I2").Value = 4584.091
$ws.Range("J2").Value = 50.81579
$ws.Range("K2").Value = 27504.546
$ws.Range("L2").Value = 304.89474
$ws.Range("M2").Value = -27391.546
$ws.Range("N2").Value = -530.89474
$ws.Range("H55").Value = 5408983
$ws.Range("J55").Value = 8339833
$ws.Range("L55").Value = 25019499
$ws.Range("N55").Value = -25019853
$ws.Range("H64").Value = 799
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 799
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H115").Value = 850.6667
$ws.Range("I115").Value = 854.2857
$ws.Range("K115").Value = 2562.8571
$ws.Range("M115").Value = -1387.8571
$ws.Range("H131").Value = 12866657
$ws.Range("J131").Value = 20842856
$ws.Range("L131").Value = 62528568
$ws.Range("N131").Value = -62538648
$ws.Range("H140").Value = 2219.1667
$ws.Range("I140").Value = 1812.25
$ws.Range("K140").Value = 5436.75
$ws.Range("M140").Value = -256.75
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 59626.75
$ws.Range("J26").Value = 65753.5
$ws.Range("L26").Value = 65753.5
$ws.Range("N26").Value = -66313.5
$ws.Range("H50").Value = 59626.75
$ws.Range("J50").Value = 65753.5
$ws.Range("L50").Value = 65753.5
$ws.Range("N50").Value = -66749.5
$ws.Range("H53").Value = 19999
$ws.Range("J53").Value = 19999
$ws.Range("L53").Value = 19999
$ws.Range("N53").Value = -21261
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()
$ws.Range("H102").Value = 3041
$ws.Range("I102").Value = 2477.9412
$ws.Range("K102").Value = 2477.9412
$ws.Range("M102").Value = -855.9412000000002
$ws.Range("H132").Value = 6597.875
$ws.Range("I132").Value = 7971.75
$ws.Range("K132").Value = 23915.25
$ws.Range("M132").Value = -21385.25
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 632
$ws.Range("I55").Value = 270.82352
$ws.Range("K55").Value = 270.82352
$ws.Range("M55").Value = -97.82351999999997
$ws.Range("H136").Value = 608733
$ws.Range("J136").Value = 2887
$ws.Range("L136").Value = 8661
$ws.Range("N136").Value = -13761
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()
$ws.Range("H126").Value = 12827993
$ws.Range("I126").Value = 20840272
$ws.Range("J126").Value = 8347.799999999999
$ws.Range("K126").Value = 62520816
$ws.Range("L126").Value = 25043.4
$ws.Range("M126").Value = -62518346
$ws.Range("N126").Value = -29983.4
$ws.Range("H132").Value = 2018.0555
$ws.Range("I132").Value = 1911.25
$ws.Range("J132").Value = 2872.5
$ws.Range("K132").Value = 5733.75
$ws.Range("L132").Value = 8617.5
$ws.Range("M132").Value = -3203.75
$ws.Range("N132").Value = -13677.5
